# Update countries & provincias Spain
# - Refresh the "datos actualizados" timestamp
# - Update COVID-19 counters for a handful of countries with newer figures
# - A few countries leap-frogged their neighbour in the ranking, so the
#   country label and figures for those two rows swap places

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 7 de Agosto de 2020 a las 13:59"

# --- Straight numeric refreshes (country stays on the same row) -------
# Row 6: India
$ws.Range("B6").Value = 2033847
$ws.Range("C6").Value = 8438
$ws.Range("D6").Value = 1381214
$ws.Range("E6").Value = 610948
$ws.Range("G6").Value = 47
$ws.Range("H6").Value = 41685

# Row 28: Catar
$ws.Range("B28").Value = 112383
$ws.Range("C28").Value = 291
$ws.Range("D28").Value = 109142
$ws.Range("E28").Value = 3061
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 180

# Row 43: Emiratos Arabes Unidos
$ws.Range("B43").Value = 62061
$ws.Range("C43").Value = 216
$ws.Range("D43").Value = 56015
$ws.Range("E43").Value = 5690
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 356

# Row 52: Barein
$ws.Range("E52").Value = 2786
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 158

# Row 79: Estado de Palestina
$ws.Range("B79").Value = 13722
$ws.Range("C79").Value = 324
$ws.Range("D79").Value = 7210
$ws.Range("E79").Value = 6418
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 94

# Row 80: Bosnia y Herzegovina
$ws.Range("B80").Value = 13687
$ws.Range("C80").Value = 291
$ws.Range("D80").Value = 7373
$ws.Range("E80").Value = 5920
$ws.Range("G80").Value = 10
$ws.Range("H80").Value = 394

# Row 82: Madagascar
$ws.Range("B82").Value = 12708
$ws.Range("C82").Value = 182
$ws.Range("D82").Value = 10412
$ws.Range("E82").Value = 2161
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 135

# Row 85: Senegal
$ws.Range("B85").Value = 10887
$ws.Range("C85").Value = 172
$ws.Range("D85").Value = 7186
$ws.Range("E85").Value = 3476
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 225

# Row 87: Consejo Danes para los Refugiados
$ws.Range("B87").Value = 9355
$ws.Range("C87").Value = 46
$ws.Range("D87").Value = 8174
$ws.Range("E87").Value = 963
$ws.Range("G87").Value = 3
$ws.Range("H87").Value = 218

# Row 95: Zambia
$ws.Range("B95").Value = 7486
$ws.Range("C95").Value = 322
$ws.Range("D95").Value = 6264
$ws.Range("E95").Value = 1022
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 200

# --- Ranking swaps (country + figures trade rows) ----------------------
# Rows 54/55: Ghana overtakes Armenia
$ws.Range("A54").Value = "Ghana"
$ws.Range("B54").Value = 40097
$ws.Range("C54").Value = 455
$ws.Range("D54").Value = 36638
$ws.Range("E54").Value = 3253
$ws.Range("G54").Value = 7
$ws.Range("H54").Value = 206

$ws.Range("A55").Value = "Armenia"
$ws.Range("B55").Value = 39985
$ws.Range("C55").Value = 166
$ws.Range("D55").Value = 32008
$ws.Range("E55").Value = 7200
$ws.Range("G55").Value = 5
$ws.Range("H55").Value = 777

# Rows 68/69: Nepal overtakes Austria
$ws.Range("A68").Value = "Nepal"
$ws.Range("B68").Value = 22214
$ws.Range("C68").Value = 464
$ws.Range("D68").Value = 15814
$ws.Range("E68").Value = 6330
$ws.Range("G68").Value = 5
$ws.Range("H68").Value = 70

$ws.Range("A69").Value = "Austria"
$ws.Range("B69").Value = 21837
$ws.Range("C69").Value = 141
$ws.Range("D69").Value = 19690
$ws.Range("E69").Value = 1427
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 720

# Rows 133/134: Islandia overtakes Benin
$ws.Range("A133").Value = "Islandia"
$ws.Range("B133").Value = 1952
$ws.Range("C133").Value = 22
$ws.Range("D133").Value = 1833
$ws.Range("E133").Value = 109
$ws.Range("H133").Value = 10

$ws.Range("A134").Value = "Benin"
$ws.Range("B134").Value = 1936
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 1600
$ws.Range("E134").Value = 298
$ws.Range("H134").Value = 38

# Rows 158/159: Vietnam overtakes Bahamas
$ws.Range("A158").Value = "Vietnam"
$ws.Range("B158").Value = 784
$ws.Range("C158").Value = 37
$ws.Range("D158").Value = 395
$ws.Range("E158").Value = 379
$ws.Range("H158").Value = 10

$ws.Range("A159").Value = "Bahamas"
$ws.Range("B159").Value = 761
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 91
$ws.Range("E159").Value = 656
$ws.Range("H159").Value = 14

# Rows 202/203: Santa Lucia overtakes Timor Oriental (figures tied, only
# the country names trade places)
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("B202").Value = 25
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 24
$ws.Range("E202").Value = 1
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

$ws.Range("A203").Value = "Timor Oriental"
$ws.Range("B203").Value = 25
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 24
$ws.Range("E203").Value = 1
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0
